$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: C21 changes from 1 to 0.5
$ws.Range("C21").Value2 = 0.5

# Insert a new row at position 22 (pushes the former row 22 down to row 23)
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new AGRI11 asset data
$ws.Range("A22").Value2 = "AGRItest"
$ws.Range("B22").Value2 = "AGRI11"
$ws.Range("C22").Value2 = 0.5
$ws.Range("D22").Value2 = 1
$ws.Range("E22").Value2 = "Agriculture"
$ws.Range("H22").Value2 = "Amazonas"
$ws.Range("L22").Value2 = "agriculture"
$ws.Range("M22").Value2 = "Soybean"

# Match the new row's height
$ws.Rows.Item(22).RowHeight = 17.25
